$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1 with the same style as the other headers (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J10
$values = @{
    2  = @(6, 6)
    3  = @(7, 7)
    4  = @(6, 6)
    5  = @(7, 7)
    6  = @(7, 7)
    7  = @(5, 5)
    8  = @(8, 8)
    9  = @(5, 5)
    10 = @(3, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
